# Insert two new weekly price rows at the top of the Pomelo / Vega Modelo de
# Temuco data block (rows 230-231), pushing the existing rows 230-246 down to
# 232-248. This mirrors the weekly update described in the commit message
# ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 230 (shifts everything below down
# by 2 rows, carrying the existing row formatting, e.g. the date style on
# column D, down with it).
$ws.Range("A230:A231").EntireRow.Insert()

# --- New row 230 ---
$ws.Cells.Item(230, 1).Value2 = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value2 = 44746
$ws.Cells.Item(230, 5).Value2 = 9
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value2 = 100102
$ws.Cells.Item(230, 8).Value = "Cítricos"
$ws.Cells.Item(230, 9).Value2 = 100102006
$ws.Cells.Item(230, 10).Value = "Pomelo"
$ws.Cells.Item(230, 11).Value = "Start Ruby"
$ws.Cells.Item(230, 12).Value = "Especial"
$ws.Cells.Item(230, 13).Value2 = 40
$ws.Cells.Item(230, 14).Value2 = 15000
$ws.Cells.Item(230, 15).Value2 = 15000
$ws.Cells.Item(230, 16).Value2 = 15000
$ws.Cells.Item(230, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(230, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(230, 19).Value2 = 1000
$ws.Cells.Item(230, 20).Value2 = 15

# --- New row 231 ---
$ws.Cells.Item(231, 1).Value2 = 10
$ws.Cells.Item(231, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(231, 3).Value = "La Araucanía"
$ws.Cells.Item(231, 4).Value2 = 44746
$ws.Cells.Item(231, 5).Value2 = 9
$ws.Cells.Item(231, 6).Value = "Fruta"
$ws.Cells.Item(231, 7).Value2 = 100102
$ws.Cells.Item(231, 8).Value = "Cítricos"
$ws.Cells.Item(231, 9).Value2 = 100102006
$ws.Cells.Item(231, 10).Value = "Pomelo"
$ws.Cells.Item(231, 11).Value = "Start Ruby"
$ws.Cells.Item(231, 12).Value = "Primera"
$ws.Cells.Item(231, 13).Value2 = 140
$ws.Cells.Item(231, 14).Value2 = 12000
$ws.Cells.Item(231, 15).Value2 = 12000
$ws.Cells.Item(231, 16).Value2 = 12000
$ws.Cells.Item(231, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(231, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(231, 19).Value2 = 800
$ws.Cells.Item(231, 20).Value2 = 15
